# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 2023-10-03 (serial 45202) to 2023-10-04 (serial 45203).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 211 }

$range = $ws.Range("C2:C$lastRow")
foreach ($cell in $range.Cells) {
    if ($cell.Value2 -eq 45202) {
        $cell.Value2 = 45203
    }
}
